$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.168.46'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '1.584.57'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '1.808.60'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '1.579.19'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '26.170.73'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.84'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  -3.03%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('E30').Value = '  -2.01%  '
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('D33').Value = '1.405.82'
$ws.Range('E33').Value = '  +7.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.589'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('E40').Value = '  +3.98%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.941'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -15.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.766'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').Value = '1.720.16'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.91'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.37%  '
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.17%  '
